$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.012.18'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '1.679.20'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("E6").Value = '  +1.43%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0620'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.38'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.33%  '
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").Value = '1.913.77'
$ws.Range("E12").Value = '  +0.25%  '
$ws.Range("D13").Value = '1.677.12'
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("E15").Value = '  +1.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.68%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '27.019.99'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '236.68'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").Value = '0.0₃0735'
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("E24").Value = '  -2.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.113'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.37%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").Value = '1.483.05'
$ws.Range("E33").Value = '  +2.12%  '
$ws.Range("E34").Value = '  +1.07%  '
$ws.Range("E35").Value = '  +5.22%  '
$ws.Range("E36").Value = '  +0.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.583'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.61%  '
$ws.Range("E38").Value = '  +2.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.907'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = '  +1.10%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("E43").Value = '  +1.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '67.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.75%  '
$ws.Range("D45").Value = '1.820.25'
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("E48").Value = '  +2.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.53'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.46%  '
$ws.Range("E50").Value = '  +1.59%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0508'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.08%  '
